# BOM.xlsx update: fixed design rule violations and added input diode
#
# - Row 7  (item 1,  CN header): designator list corrected (CN1,CN2,CN3 -> CN2,CN3,CN15)
# - Row 8  (item 2,  CN4):  part number + description/footprint corrected to a 3-circuit
#                           vertical header (2POS -> 3POS)
# - Row 9  (item 3,  CN5):  part number + description/footprint corrected to a 10-position
#                           right-angle header (8 circuits/2x4POS -> 10POS/2x5POS)
# - Row 10 (item 4,  CN6-9): designator list extended with CN13, CN14; qty 4 -> 6
# - Row 22 (item 16, new):  input protection diode D1 added

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Item 1 (row 7): designator list fix ---
$ws.Range("B7").Value = "CN2,CN3,CN15"

# --- Item 2 (row 8): CN4 header corrected to 3-position ---
$ws.Range("E8").Value = 1719710003
$ws.Range("F8").Value = "CONN HEADER VERT 3POS 2.54MM"
$ws.Range("G8").Value = "1x3POS"

# --- Item 3 (row 9): CN5 header corrected to 10-position right-angle ---
$ws.Range("E9").Value = 768250010
$ws.Range("F9").Value = "CONN HEADER R/A 10POS 5.7MM"
$ws.Range("G9").Value = "2x5POS"

# --- Item 4 (row 10): designators extended, qty corrected ---
$ws.Range("B10").Value = "CN6,CN7,CN8,CN9,CN13,CN14"
$ws.Range("C10").Value = 6

# --- Item 16 (row 22, new): input diode D1 ---
$ws.Range("A22").Value = 16
$ws.Range("B22").Value = "D1"
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = "Diodes Incorporated"
$ws.Range("E22").Value = "S1A-13-F"
$ws.Range("F22").Value = "DIODE GEN PURP 50V 1A SMA"
$ws.Range("G22").Value = "DO-214AC"
$ws.Range("H22").Value = "SMD"

# --- cosmetic: widen designator column B, move active selection ---
$ws.Columns.Item(2).ColumnWidth = 26
$ws.Range("G17").Select()
